# Kane Meeting 001 (1-10-15).docx - proofing-mark pass
#
# Word's spell/grammar checker has just been run over the document:
#   - "Quantbot" was corrected to "Quantmod" (the real R package name).
#   - A handful of words/phrases Word doesn't recognise (xtf, xbrl, MSCI,
#     e.t.c, the "a ?" / "splits." construction) got wrapped in the
#     w:proofErr spell/grammar squiggle markers that Word leaves in the
#     saved XML.
#   - The stray "_GoBack" bookmark (Word's "last edit position" marker)
#     moved from the final paragraph to the "xbrl" bullet, since that is
#     now the last place text was actually touched.
#
# The Word object model has no direct "add a proofErr" API (those marks
# are a side effect of live spell-check, not something exposed as a
# settable property), so each affected paragraph is rewritten wholesale
# via Range.InsertXML with the exact run/proofErr/bookmark structure we
# want. Replacing a paragraph's *full* Range (Paragraphs(n).Range, which
# includes the end-of-paragraph mark) with a complete <w:p> is what makes
# InsertXML behave as a true replace rather than an insert-before.

$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphXml {
    param(
        [int]$Index,
        [string]$InnerXml
    )
    $range = $d.Paragraphs($Index).Range
    $range.InsertXML("<w:p $wns>$InnerXml</w:p>")
}

# --- Paragraph 1: "(Lines or words ending with a ? indicate ...)" -------
# Grammar-flag the "a ?" construction.
$xml1 = '<w:r><w:t xml:space="preserve">(Lines or words ending with </w:t></w:r>'
$xml1 += '<w:proofErr w:type="gramStart"/>'
$xml1 += '<w:r><w:t>a ?</w:t></w:r>'
$xml1 += '<w:proofErr w:type="gramEnd"/>'
$xml1 += '<w:r><w:t xml:space="preserve"> indicate incomplete confidence in transcription.)</w:t></w:r>'
Set-ParagraphXml 1 $xml1

# --- "Quantbot" -> "Quantmod" (R tools and resources list) --------------
$xml9 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'
$xml9 += '<w:proofErr w:type="spellStart"/>'
$xml9 += '<w:r><w:t>Quantmod</w:t></w:r>'
$xml9 += '<w:proofErr w:type="spellEnd"/>'
Set-ParagraphXml 9 $xml9

# --- "xtf (?)" ------------------------------------------------------------
$xml10 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'
$xml10 += '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>'
$xml10 += '<w:r><w:t>xtf</w:t></w:r>'
$xml10 += '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>'
$xml10 += '<w:r><w:t xml:space="preserve"> (?)</w:t></w:r>'
Set-ParagraphXml 10 $xml10

# --- "xbrl" - gains the spell-check wrap *and* the _GoBack bookmark -------
$xml11 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'
$xml11 += '<w:proofErr w:type="spellStart"/>'
$xml11 += '<w:r><w:t>xbrl</w:t></w:r>'
$xml11 += '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$xml11 += '<w:proofErr w:type="spellEnd"/>'
Set-ParagraphXml 11 $xml11

# --- "MSCI" (Resources for Financial Information list) -------------------
$xml17 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>'
$xml17 += '<w:proofErr w:type="spellStart"/>'
$xml17 += '<w:r><w:t>MSCI</w:t></w:r>'
$xml17 += '<w:proofErr w:type="spellEnd"/>'
Set-ParagraphXml 17 $xml17

# --- Final bullet: "How various stock/e.t.c. listing services ..." -------
# "e.t.c" gets spell-flagged; the old _GoBack bookmark is dropped here
# (it moved to the "xbrl" bullet above); "splits." gets grammar-flagged.
$xml27 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>'
$xml27 += '<w:r><w:t>How various stock/</w:t></w:r>'
$xml27 += '<w:proofErr w:type="spellStart"/>'
$xml27 += '<w:r><w:t>e.t.c</w:t></w:r>'
$xml27 += '<w:proofErr w:type="spellEnd"/>'
$xml27 += '<w:r><w:t xml:space="preserve">. listing services </w:t></w:r>'
$xml27 += '<w:r><w:t xml:space="preserve">that we use </w:t></w:r>'
$xml27 += '<w:r><w:t xml:space="preserve">deal with things such as stock </w:t></w:r>'
$xml27 += '<w:proofErr w:type="gramStart"/>'
$xml27 += '<w:r><w:t>splits.</w:t></w:r>'
$xml27 += '<w:proofErr w:type="gramEnd"/>'
$xml27 += '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
Set-ParagraphXml 27 $xml27

Write-Output "done"
